$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 88: Md Farhan (CSE), with a live hyperlink in column D, like row 74 ---
$ws.Range("A74:F74").Copy()
$ws.Range("A88:F88").PasteSpecial(-4122)

$ws.Range("A88").Value = 45384.65997685185
$ws.Range("B88").Value = "Md Farhan"
$ws.Range("C88").Value = "B23155"
$ws.Range("D88").Value = "https://www.beecrowd.com.br/judge/en/profile/949181"
$ws.Range("E88").Value = "CSE"
$ws.Range("F88").Value = 0
$ws.Rows(88).RowHeight = 29.4

$ws.Hyperlinks.Add($ws.Range("D88"), "https://www.beecrowd.com.br/judge/en/profile/949181")

# restore the original (non-hyperlink-Add) cell formatting/style for D88
$ws.Range("D74").Copy()
$ws.Range("D88").PasteSpecial(-4122)

# --- Row 89: Saatvik Pareek (MnC), plain text in column D, like row 85 ---
$ws.Range("A85:F85").Copy()
$ws.Range("A89:F89").PasteSpecial(-4122)

$ws.Range("A89").Value = 45384.660497685189
$ws.Range("B89").Value = "Saatvik Pareek"
$ws.Range("C89").Value = "B23411"
$ws.Range("D89").Value = "https://www.beecrowd.com.br/judge/en/profile/949144"
$ws.Range("E89").Value = "MnC"
$ws.Range("F89").Value = 0
$ws.Rows(89).RowHeight = 27.6

# --- Update the view so it matches where the workbook was left scrolled/selected ---
$ws.Range("F96").Select() | Out-Null
